$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.554153919219971
$ws.Range("B1").Value = 1.773665547370911
$ws.Range("C1").Value = 1.846433162689209
$ws.Range("D1").Value = 2.23274302482605
$ws.Range("E1").Value = 3.176971197128296
